# Applies the "Unraveling the Mysteries of Space" -> "Exploring the
# Convergence of Art and Science" rewrite to $word.ActiveDocument.
#
# Strategy: walk the document in order with a forward-only cursor,
# using Range.Find.Execute to locate each text fragment relative to
# that cursor (so repeated fragments, e.g. the many lone "." runs,
# always resolve to the correct occurrence). Plain text fragments are
# replaced in place via Range.Text, which keeps them inside their
# original <w:r> (preserving rPr/run boundaries). Where the edit needs
# a brand-new run (e.g. splitting one sentence into two), we toggle
# Font.Bold on/off on the freshly written text: this is a no-op
# formatting-wise but forces the engine to keep that text in its own
# run instead of silently re-merging it with a same-formatted neighbor.

$d = $word.ActiveDocument

$script:cursor = 0

function Protect-Range($rng) {
    # Force a run boundary around $rng without changing its visible
    # formatting (set a property, then set it back).
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}

function Find-Forward($text) {
    $rng = $d.Range($script:cursor, $d.Content.End)
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find-Forward: text not found after cursor $($script:cursor): $text"
    }
    return $rng
}

# Replace $oldText (search forward from cursor) with $newText, keep it
# as its own run, and advance the cursor past it.
function Replace-Run($oldText, $newText) {
    $rng = Find-Forward $oldText
    $start = $rng.Start
    $rng.Text = $newText
    $protectRng = $d.Range($start, $start + $newText.Length)
    Protect-Range $protectRng
    $script:cursor = $start + $newText.Length
}

# Locate $text forward from cursor, leave it untouched but make sure
# it keeps its own run, then advance the cursor past it.
function Protect-Text($text) {
    $rng = Find-Forward $text
    Protect-Range $rng
    $script:cursor = $rng.End
}

# Insert brand-new text right after the cursor position, give it its
# own run, and advance the cursor past it.
function Insert-NewRun($newText) {
    $start = $script:cursor
    $rng = $d.Range($start, $start)
    $rng.InsertAfter($newText)
    $protectRng = $d.Range($start, $start + $newText.Length)
    Protect-Range $protectRng
    $script:cursor = $start + $newText.Length
}

# Move the cursor past the next occurrence of $text without altering
# anything (used to skip over <w:br/> runs etc. that Find can match
# via their surrounding text, e.g. nothing needed here, kept for
# clarity/extension).
function Skip-Forward($text) {
    $rng = Find-Forward $text
    $script:cursor = $rng.End
}

# ---------------------------------------------------------------
# Title
# ---------------------------------------------------------------
Replace-Run "Unraveling the Mysteries of Space: A Journey Through the Galaxy" "Exploring the Convergence of Art and Science"

# ---------------------------------------------------------------
# Byline: "Dr" + "." + " Emily Carter"  ->  single run "Carissa Fernandez"
# ---------------------------------------------------------------
Replace-Run "Dr" "Carissa Fernandez"
# absorb the now-orphaned ". Emily Carter" text into the run we just wrote
$rng = Find-Forward "."
$rng2 = $d.Range($rng.Start, $rng.Start + 1 + " Emily Carter".Length)
$rng2.Text = ""
$script:cursor = $rng.Start

# ---------------------------------------------------------------
# Email
# ---------------------------------------------------------------
Replace-Run "EmilyCarterPhD@cosmosresearch" "fernandezcarissa07@gmail"
Protect-Text "."
Replace-Run "edu" "com"

# ---------------------------------------------------------------
# Body paragraph 1 (first block of three sentences + blank lines)
# ---------------------------------------------------------------
Replace-Run "The vast expanse of the cosmos has captivated humanity for centuries, inspiring awe and wonder" "Art and science, often perceived as disparate disciplines, share an intrinsic connection that weaves together creativity and rationality"
Protect-Text "."
Replace-Run " From the earliest astronomers gazing up at the night sky to the modern era of space exploration, we have embarked on an ongoing journey to understand the intricacies of the universe" " They both stem from a profound curiosity and an insatiable desire to understand and express the intricacies of the world around us"
Protect-Text "."
Replace-Run " This exploration has led to profound insights into the nature of our place in the cosmos and the fundamental laws that govern the universe" " In the realm of art, we find emotions and imagination taking center stage, while in the domain of science, logic and reason lead the way"
Insert-NewRun "."
Insert-NewRun " However, upon closer examination, the boundaries between these two seemingly contrasting realms begin to blur, revealing a captivating interplay that has shaped human thought and culture throughout history"
Protect-Text "."

Replace-Run "We have witnessed the birth and death of stars, the formation of galaxies, and the enigmatic phenomena of black holes" "The convergence of art and science is evident in the ways that artistic expression can illuminate scientific concepts, making them more accessible and engaging"
Protect-Text "."
Replace-Run " We have discovered planets orbiting distant suns, raising questions about the potential for life beyond Earth" " Through paintings, sculptures, music, and literature, artists have the unique ability to translate complex scientific phenomena into forms that resonate with our senses and emotions, fostering a deeper understanding and appreciation of the natural world"
Protect-Text "."
Replace-Run " The mysteries of space continue to beckon us, fueling our insatiable curiosity and driving our quest for knowledge" " Conversely, science provides art with a rich tapestry of inspiration, offering artists a boundless realm of forms, colors, and patterns to draw upon"
Insert-NewRun "."
Insert-NewRun " The breathtaking beauty of a starry night sky, the intricate structure of a flower, or the rhythmic pulse of a heartbeat can all serve as muses, igniting the creative spark in an artist's mind"
Protect-Text "."

Replace-Run "With each new discovery, we deepen our understanding of the universe and our place within it" "Furthermore, both art and science share a common goal: to communicate ideas and inspire thought"
Protect-Text "."
Replace-Run " We unravel the secrets of cosmic evolution, unraveling the history of the universe and tracing its trajectory into the future" " Artists strive to convey their perspectives, emotions, and experiences through their works, while scientists aim to share their findings and insights with the world"
Protect-Text "."
Replace-Run " Our journey through the galaxy is a testament to humanity's enduring quest for knowledge and our unwavering fascination with the boundless mysteries of space" " Both disciplines rely on effective communication to engage their audiences, whether it be through the evocative power of imagery or the persuasive force of logical argument"
Insert-NewRun "."
Insert-NewRun " The ability to effectively communicate complex concepts is essential for both artists and scientists, and it is through this shared purpose that they find common ground"
Protect-Text "."

# ---------------------------------------------------------------
# Summary heading (unchanged) then Summary paragraph
# ---------------------------------------------------------------
Skip-Forward "Summary"

Replace-Run "Our exploration of space has yielded remarkable insights into the nature and history of the universe" "The convergence of art and science is a testament to the multifaceted nature of human understanding"
Protect-Text "."
Replace-Run " We have witnessed celestial wonders, from star formations to black holes, and discovered planets beyond our solar system" " Through their unique perspectives, artists and scientists complement each other, offering a holistic approach to comprehending the universe"
Protect-Text "."
Replace-Run " The pursuit of space exploration continues to drive our quest for knowledge, captivating humanity with its profound implications for our understanding of the universe and our place within it" " Art illuminates the emotional and intuitive dimensions of existence, while science provides a framework for rational inquiry and empirical evidence"
Insert-NewRun "."
Insert-NewRun " Together, they create a dynamic interplay that enriches our understanding of the world and fuels the progress of human thought and culture"
Protect-Text "."

# ---------------------------------------------------------------
# New trailing empty paragraph at the very end of the document body.
# ---------------------------------------------------------------
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRng.InsertParagraphAfter()
